$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# --- Row 8 content changes -------------------------------------------------
# "Fecha de término del periodo que se informa": the report now covers the
# whole fiscal year (annual) instead of just the first quarter, so the end
# date moves from 2021-03-31 to 2021-12-31.
$ws.Range("C8").Value = 44561

# "Hipervínculo a la Cuenta Pública consolidada" (D8) keeps pointing at the
# same published URL/hyperlink - no content change needed there.

# "Nota" (H8) now documents why the 2020 public account hyperlink is being
# referenced even though it was actually delivered in March.
$ws.Range("H8").Value = "La cuenta Públicas el ejercicio fiscal 2020, se entrega en la CACEH En marzo , motivo por el cual se indica el Hipervínculo "

# Row grows taller to accommodate the wrapped note text.
$ws.Rows.Item(8).RowHeight = 60

# --- Column width tweaks -----------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 98.7109375
$ws.Columns.Item(5).ColumnWidth = 73.140625

# --- View / selection state ------------------------------------------------
$ws.Activate()
$ws.Range("C8").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
